$d = $word.ActiveDocument

# --- Change 1: "W1 > W3" -> "W1 or W5 > W3" ---------------------------------
# The original paragraph has three runs: "W1 ", ">", " W3". A brand-new run
# containing "or W5 " (28pt, matching the surrounding text) must be spliced
# in between the "W1 " run and the ">" run, WITHOUT merging into the "W1 "
# run (the diff shows them as two separate <w:r> elements even though both
# end up with identical <w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>).
#
# Plain Range.InsertAfter()/Font assignment normalizes (merges) adjacent runs
# that resolve to identical formatting, so instead we delete the existing
# "W1 " run and re-insert it together with the new "or W5 " run as raw OOXML
# via Range.InsertXML at the paragraph's start (InsertXML splices inline,
# without forcing a run-merge, when targeting the exact start of a
# paragraph).
$rng = $d.Content
$rng.Find.Execute("W1 ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraStart = $rng.Start
$rng.Delete()

$insertRng = $d.Range($paraStart, $paraStart)
$openXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r w:rsidRPr="00F70B39"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">W1 </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">or W5 </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertRng.InsertXML($openXml)

# --- Change 2: "Too Weirdly Straight" -> "Too Straight" ---------------------
$d.Content.Find.Execute("Too Weirdly Straight", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Too Straight", 2)
